# Update the "取得日時" (retrieved-at) timestamp in column A for rows 2-11
# on the "ランサーズ" sheet from 2025-11-16 18:22:52 to 2025-11-16 18:30:31.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

for ($row = 2; $row -le 11; $row++) {
    $ws.Cells.Item($row, 1).Value = "2025-11-16 18:30:31"
}
